# Update column G ("K") values for rows 2-18 to reflect the new
# strikeout/K counts computed when the save_data was regenerated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 7
    3  = 2
    4  = 6
    5  = 3
    6  = 2
    7  = 3
    8  = 11
    9  = 8
    10 = 7
    11 = 11
    12 = 3
    13 = 4
    14 = 5
    15 = 7
    16 = 0
    17 = 5
    18 = 5
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
